# Update countries & provincias Spain
# Applies the 21-Jun-2020 15:13 -> 16:30 data refresh:
#  - timestamp footer cell
#  - a handful of country rows whose relative ranking shuffled (so the
#    country name shown in column A for that row changed)
#  - refreshed case/death counters for the affected rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Footer timestamp (A1) ----
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 16:30"

# ---- Country name re-ordering (same row, new country) ----
$ws.Range("A80").Value  = "Republica de Macedonia"
$ws.Range("A81").Value  = "Haiti"

$ws.Range("A115").Value = "Libano"
$ws.Range("A116").Value = "Eslovaquia"
$ws.Range("A117").Value = "Guinea-Bisau"

$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("A209").Value = "Santa Sede"

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("A214").Value = "Papua Nueva Guinea"

# ---- Refreshed numeric figures (Casos totales, Nuevos casos, Casos
#      activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ----

function Set-Row {
    param($ws, $row, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

Set-Row $ws 4   2332056 1478 973061 1236975 0 40 122020
Set-Row $ws 7   415624  3897 230586 171691  0 70 13347
Set-Row $ws 8   304331  1221 0      0       0 43 42632
Set-Row $ws 35  41204   0    12728  27476   0 8  1000
Set-Row $ws 40  31292   49   29000  336     0 0  1956
Set-Row $ws 57  13953   0    7896   5592    0 1  465
Set-Row $ws 59  12894   91   11947  686     0 1  261
Set-Row $ws 70  8743    1    8138   361     0 0  244
Set-Row $ws 80  5106    101  1926   2942    0 5  238
Set-Row $ws 81  5077    97   24     4965    0 1  88
Set-Row $ws 115 1587    51   1068   487     0 0  32
Set-Row $ws 116 1587    1    1447   112     0 0  28
Set-Row $ws 117 1541    0    153    1371    0 0  17
Set-Row $ws 161 290     3    200    84      0 0  6
Set-Row $ws 208 12      0    11     0       0 0  1
Set-Row $ws 209 12      0    12     0       0 0  0
Set-Row $ws 213 8       0    7      0       0 0  1
Set-Row $ws 214 8       0    8      0       0 0  0
